$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.21
$ws.Range("C2").Value = 0.5266666666666666
$ws.Range("J2").Value = 0.02666666666666667
$ws.Range("P2").Value = 0.1533333333333333
$ws.Range("S2").Value = 0.08333333333333333
$ws.Range("B3").Value = 0.01219512195121951
$ws.Range("C3").Value = 0.03658536585365853
$ws.Range("J3").Value = 0.0426829268292683
$ws.Range("P3").Value = 0.7195121951219512
$ws.Range("S3").Value = 0.1890243902439024
$ws.Range("J4").Value = 0.04651162790697674
$ws.Range("O4").Value = 0.02325581395348837
$ws.Range("P4").Value = 0.7209302325581395
$ws.Range("S4").Value = 0.2093023255813954
$ws.Range("B6").Value = 0.0975609756097561
$ws.Range("D6").Value = 0.01219512195121951
$ws.Range("F6").Value = 0.06097560975609756
$ws.Range("J6").Value = 0.225609756097561
$ws.Range("O6").Value = 0.03048780487804878
$ws.Range("Q6").Value = 0.1463414634146341
$ws.Range("R6").Value = 0.06707317073170732
$ws.Range("S6").Value = 0.3597560975609756
$ws.Range("B7").Value = 0.1655172413793103
$ws.Range("D7").Value = 0.04827586206896552
$ws.Range("F7").Value = 0.06206896551724138
$ws.Range("J7").Value = 0.1103448275862069
$ws.Range("O7").Value = 0.006896551724137931
$ws.Range("Q7").Value = 0.1724137931034483
$ws.Range("R7").Value = 0.07586206896551724
$ws.Range("S7").Value = 0.3586206896551724
$ws.Range("B8").Value = 0.1206896551724138
$ws.Range("D8").Value = 0.01149425287356322
$ws.Range("F8").Value = 0.04597701149425287
$ws.Range("J8").Value = 0.1293103448275862
$ws.Range("O8").Value = 0.005747126436781609
$ws.Range("Q8").Value = 0.1408045977011494
$ws.Range("R8").Value = 0.09770114942528736
$ws.Range("S8").Value = 0.4482758620689655
$ws.Range("B9").Value = 0.09289617486338798
$ws.Range("D9").Value = 0.04371584699453552
$ws.Range("E9").Value = 0.00546448087431694
$ws.Range("F9").Value = 0.07103825136612021
$ws.Range("J9").Value = 0.08196721311475409
$ws.Range("O9").Value = 0.01639344262295082
$ws.Range("Q9").Value = 0.180327868852459
$ws.Range("R9").Value = 0.0546448087431694
$ws.Range("S9").Value = 0.453551912568306
$ws.Range("B10").Value = 0.1281800391389432
$ws.Range("D10").Value = 0.02250489236790607
$ws.Range("F10").Value = 0.0675146771037182
$ws.Range("J10").Value = 0.1477495107632094
$ws.Range("O10").Value = 0.01859099804305284
$ws.Range("Q10").Value = 0.1692759295499021
$ws.Range("R10").Value = 0.07827788649706457
$ws.Range("S10").Value = 0.3679060665362035
$ws.Range("F11").Value = 0.004608294930875576
$ws.Range("G11").Value = 0.1612903225806452
$ws.Range("J11").Value = 0.06451612903225806
$ws.Range("K11").Value = 0.1935483870967742
$ws.Range("L11").Value = 0.5622119815668203
$ws.Range("S11").Value = 0.01382488479262673
$ws.Range("G12").Value = 0.7230769230769231
$ws.Range("J12").Value = 0.1769230769230769
$ws.Range("K12").Value = 0.01538461538461539
$ws.Range("L12").Value = 0.06153846153846154
$ws.Range("S12").Value = 0.02307692307692308
$ws.Range("F13").Value = 0.03448275862068965
$ws.Range("G13").Value = 0.7241379310344828
$ws.Range("J13").Value = 0.1724137931034483
$ws.Range("S13").Value = 0.06896551724137931
$ws.Range("J14").Value = 1
$ws.Range("F15").Value = 0.015625
$ws.Range("I15").Value = 0.09895833333333333
$ws.Range("J15").Value = 0.3385416666666667
$ws.Range("K15").Value = 0.0625
$ws.Range("M15").Value = 0.005208333333333333
$ws.Range("O15").Value = 0.05729166666666666
$ws.Range("S15").Value = 0.2552083333333333
$ws.Range("F16").Value = 0.005263157894736842
$ws.Range("H16").Value = 0.1894736842105263
$ws.Range("I16").Value = 0.08947368421052632
$ws.Range("J16").Value = 0.4052631578947368
$ws.Range("K16").Value = 0.1210526315789474
$ws.Range("M16").Value = 0.03684210526315789
$ws.Range("N16").Value = 0.005263157894736842
$ws.Range("O16").Value = 0.03157894736842105
$ws.Range("S16").Value = 0.1157894736842105
$ws.Range("F17").Value = 0.0231023102310231
$ws.Range("H17").Value = 0.2145214521452145
$ws.Range("I17").Value = 0.1023102310231023
$ws.Range("J17").Value = 0.3201320132013201
$ws.Range("K17").Value = 0.08580858085808581
$ws.Range("M17").Value = 0.0231023102310231
$ws.Range("O17").Value = 0.05280528052805281
$ws.Range("S17").Value = 0.1782178217821782
$ws.Range("F18").Value = 0.0352112676056338
$ws.Range("H18").Value = 0.1690140845070423
$ws.Range("I18").Value = 0.09154929577464789
$ws.Range("J18").Value = 0.4647887323943662
$ws.Range("K18").Value = 0.08450704225352113
$ws.Range("O18").Value = 0.04225352112676056
$ws.Range("S18").Value = 0.1126760563380282
$ws.Range("F19").Value = 0.01210428305400372
$ws.Range("H19").Value = 0.1815642458100559
$ws.Range("I19").Value = 0.09683426443202979
$ws.Range("J19").Value = 0.3780260707635009
$ws.Range("K19").Value = 0.09310986964618249
$ws.Range("M19").Value = 0.01675977653631285
$ws.Range("O19").Value = 0.08379888268156424
$ws.Range("S19").Value = 0.1378026070763501
